# Bug fix: connector lines were saved with a near-zero weight (w="1" EMU,
# i.e. ~0.00008pt) instead of the intended hairline width of 1pt
# (w="12700" EMU). Walk every shape on the slide and, for each straight
# connector, set its Line.Weight to 1 point so PowerPoint re-emits the
# <a:ln> element with w="12700".

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    $shapeCount = $s.Shapes.Count

    for ($i = 1; $i -le $shapeCount; $i++) {
        $shp = $s.Shapes.Item($i)

        # msoLine connector shapes (Type 9) carry the stray hairline weight.
        if ($shp.Type -eq 9) {
            if ($shp.Line.Weight -lt 1) {
                $shp.Line.Weight = 1
            }
        }
    }
}
